# Apply the "aspects_data" update:
#  - expand the description of the "engineering_material" aspect (row 13)
#  - append four new aspect rows (material, material_group, material_category, layer)
#  - widen column B slightly and move the active selection, matching the
#    author's final view state

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fix the existing "engineering_material" description (row 13, column C) ---
$ws.Range("C13").Value = "Engineering materials considered, subset of generic materials M"

# --- copy the formatting of the last existing data row (26) down onto the
#     four new rows so the new cells pick up the same styles ---
$ws.Range("A26:F26").Copy($ws.Range("A27:F30"))

# --- row 27: generic material aspect ---
$ws.Range("A27").Value = 26
$ws.Range("B27").Value = "material"
$ws.Range("C27").Value = "generic material, used in MFA and LCA to denote goods and substances"
$ws.Range("D27").Value = 5
$ws.Range("E27").Value = "M"
$ws.Range("F27").Value = "Material"

# --- row 28: material_group aspect ---
$ws.Range("A28").Value = 27
$ws.Range("B28").Value = "material_group"
$ws.Range("C28").Value = "categories of materials, such as 'reference product', 'resource, in ground', 'waste produced', used in LCI"
$ws.Range("D28").Value = 5
$ws.Range("E28").Value = "G"
$ws.Range("F28").Value = "material Group"

# --- row 29: material_category aspect ---
$ws.Range("A29").Value = 28
$ws.Range("B29").Value = "material_category"
$ws.Range("C29").Value = "broad material groups 'product', 'waste', and 'elementary', used in LCI"
$ws.Range("D29").Value = 5
$ws.Range("E29").Value = "C"
$ws.Range("F29").Value = "material Category"

# --- row 30: layer aspect ---
$ws.Range("A30").Value = 29
$ws.Range("B30").Value = "layer"
$ws.Range("C30").Value = "layer of qantification: mass, volume, energy, radioactivity, monetary, …"
$ws.Range("D30").Value = 12
$ws.Range("E30").Value = "L"
$ws.Range("F30").Value = "Layer"

# --- widen column B a bit to fit the new aspect names ---
$ws.Range("B:B").ColumnWidth = 22.44140625

# --- match the author's final selection/view state ---
$ws.Range("O24:O25").Select()
